$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "collect_region_properties.lutFile" parameter row.
# It duplicated "parcellation.lutFile" and is being dropped from the sheet;
# deleting the entire row shifts all subsequent rows up by one, which also
# removes the now-unused shared string from the workbook's string table.
$ws.Rows.Item(47).Delete()

# Update the active selection to reflect where the author left off editing.
$ws.Range("B38").Select()
